# Auto-generated edit script applying the crypto price/volume update diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.259.27"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.907.42"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5265"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3819"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07292"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9028"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08189"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.360"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008633"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "1.381.32"
$ws.Range("E18").Value = "  -27.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "27.296.49"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.084"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.527"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.742"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.836"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.833"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09258"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8289"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05079"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.231"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.363"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5820"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02003"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.081"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.545"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1524"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4927"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06198"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "

Write-Output "Applied all cell updates"
